$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must stay text (to preserve exact
# formatting such as trailing zeros / multi-dot thousand separators), so force
# a Text number format before assigning - mirrors how the source data was
# authored (plain text price strings) rather than letting Excel coerce them.
$textCells = @(
    "D5", "D6", "D9", "D10", "D11", "D16", "D20", "D22", "D23", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D41", "D42", "D45", "D47", "D48", "D49", "D51"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated coin data (price + 1h volume change columns, plus the
# couple of rows whose coin/link order was swapped).
$ws.Cells.Item(2, 4).Value = "42.910.13"
$ws.Cells.Item(2, 5).Value = "  +0.46%  "
$ws.Cells.Item(3, 4).Value = "2.537.02"
$ws.Cells.Item(3, 5).Value = "  -0.37%  "
$ws.Cells.Item(5, 4).Value = "317.17"
$ws.Cells.Item(5, 5).Value = "  -0.55%  "
$ws.Cells.Item(6, 4).Value = "96.29"
$ws.Cells.Item(6, 5).Value = "  +1.28%  "
$ws.Cells.Item(7, 5).Value = "  -0.49%  "
$ws.Cells.Item(8, 5).Value = "  +0.02%  "
$ws.Cells.Item(9, 4).Value = "0.541"
$ws.Cells.Item(9, 5).Value = "  +1.38%  "
$ws.Cells.Item(10, 4).Value = "35.64"
$ws.Cells.Item(10, 5).Value = "  -2.19%  "
$ws.Cells.Item(11, 4).Value = "0.0814"
$ws.Cells.Item(11, 5).Value = "  +0.00%  "
$ws.Cells.Item(12, 5).Value = "  -1.25%  "
$ws.Cells.Item(13, 5).Value = "  -4.50%  "
$ws.Cells.Item(14, 4).Value = "2.929.07"
$ws.Cells.Item(14, 5).Value = "  -0.22%  "
$ws.Cells.Item(15, 4).Value = "2.536.95"
$ws.Cells.Item(15, 5).Value = "  +1.17%  "
$ws.Cells.Item(16, 4).Value = "14.99"
$ws.Cells.Item(16, 5).Value = "  -6.28%  "
$ws.Cells.Item(17, 5).Value = "  -2.25%  "
$ws.Cells.Item(18, 4).Value = "42.979.35"
$ws.Cells.Item(18, 5).Value = "  +0.56%  "
$ws.Cells.Item(19, 5).Value = "  +2.35%  "
$ws.Cells.Item(20, 4).Value = "12.57"
$ws.Cells.Item(20, 5).Value = "  -3.88%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0965"
$ws.Cells.Item(21, 5).Value = "  -0.53%  "
$ws.Cells.Item(22, 4).Value = "69.69"
$ws.Cells.Item(22, 5).Value = "  -2.02%  "
$ws.Cells.Item(23, 4).Value = "253.29"
$ws.Cells.Item(23, 5).Value = "  +0.23%  "
$ws.Cells.Item(24, 5).Value = "  -0.91%  "
$ws.Cells.Item(25, 5).Value = "  +1.66%  "
$ws.Cells.Item(26, 4).Value = "27.01"
$ws.Cells.Item(26, 5).Value = "  -0.99%  "
$ws.Cells.Item(27, 4).Value = "0.999"
$ws.Cells.Item(27, 5).Value = "  -0.16%  "
$ws.Cells.Item(28, 5).Value = "  +1.91%  "
$ws.Cells.Item(29, 4).Value = "40.60"
$ws.Cells.Item(29, 5).Value = "  +1.82%  "
$ws.Cells.Item(30, 4).Value = "10.35"
$ws.Cells.Item(30, 5).Value = "  +0.60%  "
$ws.Cells.Item(31, 4).Value = "5.89"
$ws.Cells.Item(31, 5).Value = "  -1.36%  "
$ws.Cells.Item(32, 4).Value = "155.97"
$ws.Cells.Item(32, 5).Value = "  -0.16%  "
$ws.Cells.Item(33, 4).Value = "19.46"
$ws.Cells.Item(33, 5).Value = "  +0.85%  "
$ws.Cells.Item(34, 4).Value = "3.37"
$ws.Cells.Item(34, 5).Value = "  +0.10%  "
$ws.Cells.Item(35, 4).Value = "2.71"
$ws.Cells.Item(35, 5).Value = "  +3.00%  "
$ws.Cells.Item(36, 2).Value = "Hedera"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(36, 4).Value = "0.0798"
$ws.Cells.Item(36, 5).Value = "  +0.65%  "
$ws.Cells.Item(37, 2).Value = "ARBITRUM"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(37, 4).Value = "2.10"
$ws.Cells.Item(37, 5).Value = "  -2.26%  "
$ws.Cells.Item(38, 5).Value = "  +1.57%  "
$ws.Cells.Item(39, 4).Value = "2.45"
$ws.Cells.Item(39, 5).Value = "  -0.26%  "
$ws.Cells.Item(40, 5).Value = "  -0.74%  "
$ws.Cells.Item(41, 4).Value = "21.83"
$ws.Cells.Item(41, 5).Value = "  -7.79%  "
$ws.Cells.Item(42, 4).Value = "3.81"
$ws.Cells.Item(42, 5).Value = "  -1.14%  "
$ws.Cells.Item(44, 5).Value = "  +0.01%  "
$ws.Cells.Item(45, 4).Value = "3.27"
$ws.Cells.Item(45, 5).Value = "  -2.64%  "
$ws.Cells.Item(46, 4).Value = "2.001.21"
$ws.Cells.Item(46, 5).Value = "  -1.43%  "
$ws.Cells.Item(47, 4).Value = "9.06"
$ws.Cells.Item(47, 5).Value = "  +1.29%  "
$ws.Cells.Item(48, 4).Value = "84.63"
$ws.Cells.Item(48, 5).Value = "  +0.12%  "
$ws.Cells.Item(49, 2).Value = "ordi"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Cells.Item(49, 4).Value = "74.88"
$ws.Cells.Item(49, 5).Value = "  +1.18%  "
$ws.Cells.Item(50, 2).Value = "RocketPoolETH"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(50, 4).Value = "2.783.82"
$ws.Cells.Item(50, 5).Value = "  -0.21%  "
$ws.Cells.Item(51, 4).Value = "104.74"
$ws.Cells.Item(51, 5).Value = "  +2.77%  "
